$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Pupil Class" column header becomes just "Class" (import class function)
$ws.Range("C1").Value = "Class"

# Replace placeholder parent-name sample data with real-looking names
$ws.Range("D2").Value = "John"
$ws.Range("E2").Value = "Mason"
$ws.Range("D3").Value = "Cater"
$ws.Range("E3").Value = "Jackson"

# Move the active selection to F3
$ws.Range("F3").Select()
